# Update "想去人数" (want-to-go count) figures on both the "展览" sheet
# and the "全部类型" sheet, which carry duplicate data.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F2").Value = 8761
    $ws.Range("F4").Value = 418
    $ws.Range("F5").Value = 144
}
